# Applies the "RestConfiguration / JavaGeneratorMicroservices" update to the
# petclinic graph workbook:
#  - nodes sheet: inserts a new "OwnerRepository.save(Owner)" method node,
#    adds three new "SpringBootApplication" Annotation nodes (one per
#    microservice), and clears the SubType of the rows that used to encode
#    that information inline (repository / SpringBootApplication -> "-").
#  - Edges sheet: adds a header row (Src/Dst/Type/Label), three new
#    "Has Annotation" edges (PetClinicApplication -> SpringBootApplication)
#    and one new "Has Method" edge (OwnerRepository -> save(Owner)).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("nodes")
$ws2 = $wb.Worksheets.Item("Edges")

# ---------------------------------------------------------------------
# 1. nodes sheet: insert the new "save(Owner)" row right before the old
#    row 36 (OwnerRepository.findById) so it lands at row 36 and pushes
#    everything else down by one.
# ---------------------------------------------------------------------
$ws1.Rows(36).Insert()

$ws1.Cells.Item(36, 1).Value = "owner.org.springframework.samples.petclinic.owner.OwnerRepository.save(Owner)"
$ws1.Cells.Item(36, 2).Value = "org.springframework.samples.petclinic.owner"
$ws1.Cells.Item(36, 3).Value = "OwnerRepository.save(Owner)"
$ws1.Cells.Item(36, 4).Value = "label"
$ws1.Cells.Item(36, 5).Value = "method"
$ws1.Cells.Item(36, 6).Value = "-"
$ws1.Cells.Item(36, 7).Value = "owner"

# The rows that used to carry "repository" / "SpringBootApplication" in the
# SubType column now just hold "-" there (those facts moved to their own
# Annotation nodes / Has Annotation edges).
$ws1.Cells.Item(37, 6).Value = "-"   # OwnerRepository.findById(Integer)
$ws1.Cells.Item(38, 6).Value = "-"   # pet1 PetClinicApplication
$ws1.Cells.Item(39, 6).Value = "-"   # pet2 PetClinicApplication
$ws1.Cells.Item(40, 6).Value = "-"   # owner PetClinicApplication

# New "SpringBootApplication" Annotation nodes, appended after row 45.
$ws1.Cells.Item(46, 1).Value = "owner.org.springframework.boot.autoconfigure.SpringBootApplication"
$ws1.Cells.Item(46, 2).Value = "org.springframework.boot.autoconfigure"
$ws1.Cells.Item(46, 3).Value = "SpringBootApplication"
$ws1.Cells.Item(46, 4).Value = "label"
$ws1.Cells.Item(46, 5).Value = "Annotation"
$ws1.Cells.Item(46, 6).Value = "-"
$ws1.Cells.Item(46, 7).Value = "owner"

$ws1.Cells.Item(47, 1).Value = "pet1.org.springframework.boot.autoconfigure.SpringBootApplication"
$ws1.Cells.Item(47, 2).Value = "org.springframework.boot.autoconfigure"
$ws1.Cells.Item(47, 3).Value = "SpringBootApplication"
$ws1.Cells.Item(47, 4).Value = "label"
$ws1.Cells.Item(47, 5).Value = "Annotation"
$ws1.Cells.Item(47, 6).Value = "-"
$ws1.Cells.Item(47, 7).Value = "pet1"

$ws1.Cells.Item(48, 1).Value = "pet2.org.springframework.boot.autoconfigure.SpringBootApplication"
$ws1.Cells.Item(48, 2).Value = "org.springframework.boot.autoconfigure"
$ws1.Cells.Item(48, 3).Value = "SpringBootApplication"
$ws1.Cells.Item(48, 4).Value = "label"
$ws1.Cells.Item(48, 5).Value = "Annotation"
$ws1.Cells.Item(48, 6).Value = "-"
$ws1.Cells.Item(48, 7).Value = "pet2"

# ---------------------------------------------------------------------
# 2. Edges sheet: add a header row, plus the new annotation/method edges.
# ---------------------------------------------------------------------
$ws2.Cells.Item(1, 1).Value = "Src"
$ws2.Cells.Item(1, 2).Value = "Dst"
$ws2.Cells.Item(1, 3).Value = "Type"
$ws2.Cells.Item(1, 4).Value = "Label"

$ws2.Cells.Item(33, 1).Value = "pet1.org.springframework.samples.petclinic.PetClinicApplication"
$ws2.Cells.Item(33, 2).Value = "pet1.org.springframework.boot.autoconfigure.SpringBootApplication"
$ws2.Cells.Item(33, 3).Value = "Has Annotation"
$ws2.Cells.Item(33, 4).Value = "label"

$ws2.Cells.Item(34, 1).Value = "pet2.org.springframework.samples.petclinic.PetClinicApplication"
$ws2.Cells.Item(34, 2).Value = "pet2.org.springframework.boot.autoconfigure.SpringBootApplication"
$ws2.Cells.Item(34, 3).Value = "Has Annotation"
$ws2.Cells.Item(34, 4).Value = "label"

$ws2.Cells.Item(35, 1).Value = "owner.org.springframework.samples.petclinic.PetClinicApplication"
$ws2.Cells.Item(35, 2).Value = "owner.org.springframework.boot.autoconfigure.SpringBootApplication"
$ws2.Cells.Item(35, 3).Value = "Has Annotation"
$ws2.Cells.Item(35, 4).Value = "label"

$ws2.Cells.Item(36, 1).Value = "owner.org.springframework.samples.petclinic.owner.OwnerRepository"
$ws2.Cells.Item(36, 2).Value = "owner.org.springframework.samples.petclinic.owner.OwnerRepository.save(Owner)"
$ws2.Cells.Item(36, 3).Value = "Has Method"
$ws2.Cells.Item(36, 4).Value = "label"

# ---------------------------------------------------------------------
# 3. View state: Edges becomes the active sheet/tab, with new selections.
# ---------------------------------------------------------------------
$ws1.Range("A45").Select()
$ws2.Activate()
$ws2.Range("A31").Select()
